$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (Spanish labels -> short English column names) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case municipality / state names in columns A and B ---
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B24").Value = "Amatenango De La Frontera"
$ws.Range("B26").Value = "Bejucal De Ocampo"
$ws.Range("B34").Value = "Comitán De Domínguez"
$ws.Range("B53").Value = "Mazapa De Madero"
$ws.Range("B62").Value = "Salto De Agua"
$ws.Range("B63").Value = "San Cristóbal De Las Casas"
$ws.Range("B106").Value = "Villa De Álvarez"
$ws.Range("A108").Value = "Ciudad De México"
$ws.Range("B112").Value = "Cuajimalpa De Morelos"
$ws.Range("B127").Value = "Coneto De Comonfort"
$ws.Range("B135").Value = "Nombre De Dios"
$ws.Range("B139").Value = "San Juan Del Río"
$ws.Range("B143").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B146").Value = "Almoloya De Alquisiras"
$ws.Range("B147").Value = "Almoloya Del Río"
$ws.Range("B150").Value = "Atizapán De Zaragoza"
$ws.Range("B154").Value = "Chapa De Mota"
$ws.Range("B162").Value = "Ecatepec De Morelos"
$ws.Range("B167").Value = "Ixtapan De La Sal"
$ws.Range("B176").Value = "Naucalpan De Juárez"
$ws.Range("B185").Value = "San Antonio La Isla"
$ws.Range("B186").Value = "San Felipe Del Progreso"
$ws.Range("B187").Value = "San Martín De Las Pirámides"
$ws.Range("B189").Value = "Soyaniquilpan De Juárez"
$ws.Range("B196").Value = "Tenango Del Valle"
$ws.Range("B203").Value = "Tlalnepantla De Baz"
$ws.Range("B208").Value = "Valle De Bravo"
$ws.Range("B209").Value = "Valle De Chalco Solidaridad"
$ws.Range("B210").Value = "Villa Del Carbón"
$ws.Range("B220").Value = "San Miguel De Allende"
$ws.Range("B221").Value = "Apaseo El Alto"
$ws.Range("B222").Value = "Apaseo El Grande"
$ws.Range("B229").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B233").Value = "Jaral Del Progreso"
$ws.Range("B240").Value = "Purísima Del Rincón"
$ws.Range("B244").Value = "San Diego De La Unión"
$ws.Range("B246").Value = "San Francisco Del Rincón"
$ws.Range("B248").Value = "San Luis De La Paz"
$ws.Range("B249").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B250").Value = "Silao De La Victoria"
$ws.Range("B253").Value = "Valle De Santiago"
$ws.Range("B259").Value = "Acapulco De Juárez"
$ws.Range("B261").Value = "Ajuchitlán Del Progreso"
$ws.Range("B262").Value = "Alcozauca De Guerrero"
$ws.Range("B264").Value = "Atenango Del Río"
$ws.Range("B266").Value = "Atoyac De Álvarez"
$ws.Range("B267").Value = "Ayutla De Los Libres"
$ws.Range("B268").Value = "Buenavista De Cuéllar"
$ws.Range("B269").Value = "Chilapa De Álvarez"
$ws.Range("B270").Value = "Chilpancingo De Los Bravo"
$ws.Range("B271").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B275").Value = "Coyuca De Benítez"
$ws.Range("B276").Value = "Coyuca De Catalán"
$ws.Range("B279").Value = "Cuetzala Del Progreso"
$ws.Range("B280").Value = "Cutzamala De Pinzón"
$ws.Range("B284").Value = "Huitzuco De Los Figueroa"
$ws.Range("B285").Value = "Iguala De La Independencia"
$ws.Range("B289").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B302").Value = "Taxco De Alarcón"
$ws.Range("B304").Value = "Técpan De Galeana"
$ws.Range("B306").Value = "Tepecoacuilco De Trujano"
$ws.Range("B308").Value = "Tixtla De Guerrero"
$ws.Range("B311").Value = "Tlapa De Comonfort"
$ws.Range("B323").Value = "Atotonilco El Grande"
$ws.Range("B327").Value = "Huasca De Ocampo"
$ws.Range("B329").Value = "Huejutla De Reyes"
$ws.Range("B332").Value = "Jacala De Ledezma"
$ws.Range("B335").Value = "Mineral Del Chico"
$ws.Range("B336").Value = "Molango De Escamilla"
$ws.Range("B338").Value = "Omitlán De Juárez"
$ws.Range("B339").Value = "Pachuca De Soto"
$ws.Range("B347").Value = "Tepehuacán De Guerrero"
$ws.Range("B348").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B350").Value = "Tezontepec De Aldama"
$ws.Range("B353").Value = "Tula De Allende"
$ws.Range("B354").Value = "Tulancingo De Bravo"
$ws.Range("B355").Value = "Zacualtipán De Ángeles"
$ws.Range("B361").Value = "Atemajac De Brizuela"
$ws.Range("B362").Value = "Atotonilco El Alto"
$ws.Range("B363").Value = "Autlán De Navarro"
$ws.Range("B368").Value = "Concepción De Buenos Aires"
$ws.Range("B372").Value = "Encarnación De Díaz"
$ws.Range("B376").Value = "Ixtlahuacán Del Río"
$ws.Range("B378").Value = "Jilotlán De Los Dolores"
$ws.Range("B380").Value = "Lagos De Moreno"
$ws.Range("B387").Value = "San Diego De Alejandría"
$ws.Range("B388").Value = "San Juanito De Escobedo"
$ws.Range("B392").Value = "Tamazula De Gordiano"
$ws.Range("B395").Value = "Teocuitatlán De Corona"
$ws.Range("B396").Value = "Tepatitlán De Morelos"
$ws.Range("B399").Value = "Tizapán El Alto"
$ws.Range("B400").Value = "Tlajomulco De Zúñiga"
$ws.Range("B406").Value = "Unión De San Antonio"
$ws.Range("B407").Value = "Unión De Tula"
$ws.Range("B412").Value = "Zapotlán El Grande"
$ws.Range("B427").Value = "Cojumatlán De Régules"
$ws.Range("B477").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B506").Value = "Puente De Ixtla"
$ws.Range("B510").Value = "Tetela Del Volcán"
$ws.Range("B511").Value = "Tlaltizapán De Zapata"
$ws.Range("B517").Value = "Zacualpan De Amilpas"
$ws.Range("B519").Value = "Ixtlán Del Río"
$ws.Range("B525").Value = "Santa María Del Oro"
$ws.Range("B538").Value = "San Nicolás De Los Garza"
$ws.Range("B541").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B546").Value = "Chalcatongo De Hidalgo"
$ws.Range("B548").Value = "Constancia Del Rosario"
$ws.Range("B550").Value = "Cuilápam De Guerrero"
$ws.Range("B551").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B552").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B553").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B555").Value = "Ixtlán De Juárez"
$ws.Range("B556").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B566").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B567").Value = "Mixistlán De La Reforma"
$ws.Range("B571").Value = "Oaxaca De Juárez"
$ws.Range("B572").Value = "Ocotlán De Morelos"
$ws.Range("B574").Value = "Putla Villa De Guerrero"
$ws.Range("B575").Value = "Reforma De Pineda"
$ws.Range("B586").Value = "San Antonino El Alto"
$ws.Range("B589").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B606").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B611").Value = "San Juan De Los Cués"
$ws.Range("B631").Value = "San Mateo Del Mar"
$ws.Range("B637").Value = "San Miguel Del Puerto"
$ws.Range("B646").Value = "San Pablo Villa De Mitla"
$ws.Range("B650").Value = "San Pedro El Alto"
$ws.Range("B678").Value = "Santa Inés Del Monte"
$ws.Range("B686").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B699").Value = "Santiago Del Río"
$ws.Range("B713").Value = "Santo Domingo De Morelos"
$ws.Range("B722").Value = "Teotitlán De Flores Magón"
$ws.Range("B723").Value = "Teotitlán Del Valle"
$ws.Range("B725").Value = "Tlacolula De Matamoros"
$ws.Range("B726").Value = "Totontepec Villa De Morelos"
$ws.Range("B728").Value = "Villa De Etla"
$ws.Range("B729").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B730").Value = "Villa De Zaachila"
$ws.Range("B732").Value = "Villa Sola De Vega"
$ws.Range("B733").Value = "Villa Talea De Castro"
$ws.Range("B736").Value = "Zimatlán De Álvarez"
$ws.Range("B747").Value = "Chalchicomula De Sesma"
$ws.Range("B764").Value = "Huehuetlán El Chico"
$ws.Range("B768").Value = "Ixcamilpa De Guerrero"
$ws.Range("B770").Value = "Izúcar De Matamoros"
$ws.Range("B774").Value = "Los Reyes De Juárez"
$ws.Range("B776").Value = "Palmar De Bravo"
$ws.Range("B790").Value = "San Salvador El Verde"
$ws.Range("B794").Value = "Tecali De Herrera"
$ws.Range("B801").Value = "Tepexi De Rodríguez"
$ws.Range("B802").Value = "Tetela De Ocampo"
$ws.Range("B806").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B823").Value = "Amealco De Bonfil"
$ws.Range("B825").Value = "Cadereyta De Montes"
$ws.Range("B829").Value = "Jalpan De Serra"
$ws.Range("B830").Value = "Landa De Matamoros"
$ws.Range("B833").Value = "Pinal De Amoles"
$ws.Range("B836").Value = "San Juan Del Río"
$ws.Range("B845").Value = "Axtla De Terrazas"
$ws.Range("B850").Value = "Ciudad Del Maíz"
$ws.Range("B856").Value = "Mexquitic De Carmona"
$ws.Range("B861").Value = "San Ciro De Acosta"
$ws.Range("B863").Value = "Santa María Del Río"
$ws.Range("B864").Value = "Soledad De Graciano Sánchez"
$ws.Range("B868").Value = "Tanquián De Escobedo"
$ws.Range("B870").Value = "Villa De Arista"
$ws.Range("B871").Value = "Villa De Guadalupe"
$ws.Range("B872").Value = "Villa De Ramos"
$ws.Range("B873").Value = "Villa De Reyes"
$ws.Range("B928").Value = "Soto La Marina"
$ws.Range("B937").Value = "Contla De Juan Cuamatzi"
$ws.Range("B939").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B940").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B957").Value = "Boca Del Río"
$ws.Range("B960").Value = "Cazones De Herrera"
$ws.Range("B969").Value = "Cosamaloapan De Carpio"
$ws.Range("B978").Value = "Hueyapan De Ocampo"
$ws.Range("B979").Value = "Ignacio De La Llave"
$ws.Range("B981").Value = "Ixhuatlán De Madero"
$ws.Range("B982").Value = "Ixhuatlán Del Sureste"
$ws.Range("B991").Value = "Juchique De Ferrer"
$ws.Range("B995").Value = "Lerdo De Tejada"
$ws.Range("B997").Value = "Martínez De La Torre"
$ws.Range("B1009").Value = "Paso De Ovejas"
$ws.Range("B1010").Value = "Paso Del Macho"
$ws.Range("B1013").Value = "Poza Rica De Hidalgo"
$ws.Range("B1020").Value = "Sayula De Alemán"
$ws.Range("B1021").Value = "Soledad De Doblado"
$ws.Range("B1034").Value = "Tlacotepec De Mejía"
$ws.Range("B1044").Value = "Vega De Alatorre"
$ws.Range("B1050").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1051").Value = "Zozocolco De Hidalgo"
$ws.Range("B1070").Value = "Nochistlán De Mejía"
$ws.Range("B1071").Value = "Noria De Ángeles"

# --- Clean up A143 (strip stray _x000D_ escape + trailing newline, and title-case) ---
$ws.Range("A143").Value = "Estado De México"

# --- Tiny floating point refresh on two percentage cells ---
$ws.Range("D246").Value = 0.009707287933094383
$ws.Range("D449").Value = 0.009557945041816007

# --- Remove trailing footer/metadata rows (1085:1089) ---
$ws.Rows("1085:1089").Delete()

